$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data1 = New-Object 'object[,]' 24,5
$data1[0,0] = 1.02
$data1[0,1] = 1.018777190134748
$data1[0,2] = 1.025507051647289
$data1[0,3] = 1.028441947054263
$data1[0,4] = 1.035878506188936
$data1[1,0] = 1.02
$data1[1,1] = 1.019837964448458
$data1[1,2] = 1.026019325364456
$data1[1,3] = 1.029421715039732
$data1[1,4] = 1.037028821821524
$data1[2,0] = 1.02
$data1[2,1] = 1.020524907074266
$data1[2,2] = 1.026344013162027
$data1[2,3] = 1.030056559725998
$data1[2,4] = 1.037774361322411
$data1[3,0] = 1.02
$data1[3,1] = 1.020813829861862
$data1[3,2] = 1.026478881440165
$data1[3,3] = 1.030323656259277
$data1[3,4] = 1.038088075393573
$data1[4,0] = 1.02
$data1[4,1] = 1.020862349024004
$data1[4,2] = 1.026501430639359
$data1[4,3] = 1.030368515121524
$data1[4,4] = 1.038140766328743
$data1[5,0] = 1.02
$data1[5,1] = 1.02052876715525
$data1[5,2] = 1.026345821692922
$data1[5,3] = 1.030060127866673
$data1[5,4] = 1.037778552048582
$data1[6,0] = 1.02
$data1[6,1] = 1.019135569909541
$data1[6,2] = 1.025681579129091
$data1[6,3] = 1.028772884070614
$data1[6,4] = 1.036267010527895
$data1[7,0] = 1.02
$data1[7,1] = 1.01668478420329
$data1[7,2] = 1.024459342914767
$data1[7,3] = 1.026511265982619
$data1[7,4] = 1.033612736623003
$data1[8,0] = 1.02
$data1[8,1] = 1.015053735864303
$data1[8,2] = 1.023610008354521
$data1[8,3] = 1.025008019577426
$data1[8,4] = 1.03184944187669
$data1[9,0] = 1.02
$data1[9,1] = 1.014348135187405
$data1[9,2] = 1.02323410637176
$data1[9,3] = 1.024358166769749
$data1[9,4] = 1.03108739017691
$data1[10,0] = 1.02
$data1[10,1] = 1.014086141422113
$data1[10,2] = 1.023093261465912
$data1[10,3] = 1.024116942290496
$data1[10,4] = 1.0308045503853
$data1[11,0] = 1.02
$data1[11,1] = 1.01414233555925
$data1[11,2] = 1.02312352824074
$data1[11,3] = 1.024168678527527
$data1[11,4] = 1.030865210521971
$data1[12,0] = 1.02
$data1[12,1] = 1.014326476710061
$data1[12,2] = 1.023222488916466
$data1[12,3] = 1.024338223818738
$data1[12,4] = 1.03106400606181
$data1[13,0] = 1.02
$data1[13,1] = 1.014439945087112
$data1[13,2] = 1.023283300548361
$data1[13,3] = 1.024442707445325
$data1[13,4] = 1.03118651971999
$data1[14,0] = 1.02
$data1[14,1] = 1.015100577992576
$data1[14,2] = 1.023634784689578
$data1[14,3] = 1.025051170569725
$data1[14,4] = 1.031900047620276
$data1[15,0] = 1.02
$data1[15,1] = 1.015515150292192
$data1[15,2] = 1.023853086440131
$data1[15,3] = 1.025433128033761
$data1[15,4] = 1.032348017321145
$data1[16,0] = 1.02
$data1[16,1] = 1.015757026439381
$data1[16,2] = 1.02397963304905
$data1[16,3] = 1.025656019899017
$data1[16,4] = 1.032609451984648
$data1[17,0] = 1.02
$data1[17,1] = 1.015839510679601
$data1[17,2] = 1.024022648879114
$data1[17,3] = 1.025732037653505
$data1[17,4] = 1.032698618457841
$data1[18,0] = 1.02
$data1[18,1] = 1.015470664081778
$data1[18,2] = 1.023829745907683
$data1[18,3] = 1.025392137006533
$data1[18,4] = 1.032299939776394
$data1[19,0] = 1.02
$data1[19,1] = 1.014272249036662
$data1[19,2] = 1.023193381061208
$data1[19,3] = 1.024288292544409
$data1[19,4] = 1.03100545965476
$data1[20,0] = 1.02
$data1[20,1] = 1.013519323194119
$data1[20,2] = 1.022786227006253
$data1[20,3] = 1.023595186264539
$data1[20,4] = 1.030192842195708
$data1[21,0] = 1.02
$data1[21,1] = 1.013918410018228
$data1[21,2] = 1.023002733729557
$data1[21,3] = 1.023962527341578
$data1[21,4] = 1.030623505428514
$data1[22,0] = 1.02
$data1[22,1] = 1.015490765290502
$data1[22,2] = 1.023840294915894
$data1[22,3] = 1.025410658770727
$data1[22,4] = 1.03232166351239
$data1[23,0] = 1.02
$data1[23,1] = 1.017317874621407
$data1[23,2] = 1.024781422877932
$data1[23,3] = 1.027095156568681
$data1[23,4] = 1.034297833428938

$ws.Range("B2:F25").Value = $data1

$data2 = New-Object 'object[,]' 24,6
$data2[0,0] = 1.028047310659273
$data2[0,1] = 1.023983701230587
$data2[0,2] = 1.028332420894459
$data2[0,3] = 1.031258752783233
$data2[0,4] = 1.038673849366234
$data2[0,5] = 1.012006349856506
$data2[1,0] = 1.028114196249099
$data2[1,1] = 1.024680361429499
$data2[1,2] = 1.02865281339994
$data2[1,3] = 1.032045993535909
$data2[1,4] = 1.039632745861786
$data2[1,5] = 1.012244501852445
$data2[2,0] = 1.028152502439901
$data2[2,1] = 1.025131167901118
$data2[2,2] = 1.028852644648525
$data2[2,3] = 1.032555648281803
$data2[2,4] = 1.040253847108459
$data2[2,5] = 1.01239841392067
$data2[3,0] = 1.028167412619882
$data2[3,1] = 1.025320691769214
$data2[3,2] = 1.028934856412712
$data2[3,3] = 1.032769968569184
$data2[3,4] = 1.040515109076832
$data2[3,5] = 1.012463073293099
$data2[4,0] = 1.0281698460087
$data2[4,1] = 1.025352513939803
$data2[4,2] = 1.028948554567337
$data2[4,3] = 1.032805957485931
$data2[4,4] = 1.040558984963684
$data2[4,5] = 1.012473927226143
$data2[5,0] = 1.028152706365405
$data2[5,1] = 1.025133700309516
$data2[5,2] = 1.028853750236068
$data2[5,3] = 1.032558511798797
$data2[5,4] = 1.040257337510414
$data2[5,5] = 1.012399278080471
$data2[6,0] = 1.02807094286734
$data2[6,1] = 1.024219136428681
$data2[6,2] = 1.028442245416833
$data2[6,3] = 1.031524751091472
$data2[6,4] = 1.038997782203664
$data2[6,5] = 1.012086873282085
$data2[7,0] = 1.027888906352644
$data2[7,1] = 1.022607722786081
$data2[7,2] = 1.027660033183786
$data2[7,3] = 1.029705111133826
$data2[7,4] = 1.036783121677486
$data2[7,5] = 1.011534943782861
$data2[8,0] = 1.027742192411911
$data2[8,1] = 1.021533564086832
$data2[8,2] = 1.027100466945094
$data2[8,3] = 1.028493356457054
$data2[8,4] = 1.035309935747672
$data2[8,5] = 1.011166039416976
$data2[9,0] = 1.027672681232815
$data2[9,1] = 1.021068470029094
$data2[9,2] = 1.026849189472056
$data2[9,3] = 1.027968971972987
$data2[9,4] = 1.034672801649342
$data2[9,5] = 1.011006076378505
$data2[10,0] = 1.027645965106361
$data2[10,1] = 1.020895716942084
$data2[10,2] = 1.026754508480939
$data2[10,3] = 1.027774239312679
$data2[10,4] = 1.03443625652106
$data2[10,5] = 1.010946625336262
$data2[11,0] = 1.027651736342748
$data2[11,1] = 1.020932772905002
$data2[11,2] = 1.02677487869296
$data2[11,3] = 1.027816007998315
$data2[11,4] = 1.034486991058778
$data2[11,5] = 1.010959379306702
$data2[12,0] = 1.027670491146611
$data2[12,1] = 1.021054190132586
$data2[12,2] = 1.026841390534216
$data2[12,3] = 1.027952874349895
$data2[12,4] = 1.034653246413062
$data2[12,5] = 1.011001162822286
$data2[13,0] = 1.027681927861248
$data2[13,1] = 1.021128999764434
$data2[13,2] = 1.026882192518695
$data2[13,3] = 1.028037208466125
$data2[13,4] = 1.034755697038808
$data2[13,5] = 1.011026902571408
$data2[14,0] = 1.027746679784026
$data2[14,1] = 1.02156443134678
$data2[14,2] = 1.02711695452225
$data2[14,3] = 1.028528164730228
$data2[14,4] = 1.035352236350497
$data2[14,5] = 1.011176650909739
$data2[15,0] = 1.027785696793537
$data2[15,1] = 1.02183757242101
$data2[15,2] = 1.027261812832631
$data2[15,3] = 1.028836212510061
$data2[15,4] = 1.035726634468873
$data2[15,5] = 1.011270523954976
$data2[16,0] = 1.027807877329898
$data2[16,1] = 1.021996893281358
$data2[16,2] = 1.027345439387844
$data2[16,3] = 1.029015921795615
$data2[16,4] = 1.03594508838923
$data2[16,5] = 1.011325256793188
$data2[17,0] = 1.027815342305294
$data2[17,1] = 1.02205121796653
$data2[17,2] = 1.027373806722037
$data2[17,3] = 1.029077203144092
$data2[17,4] = 1.036019588056826
$data2[17,5] = 1.01134391558742
$data2[18,0] = 1.027781570346936
$data2[18,1] = 1.021808266732241
$data2[18,2] = 1.027246360540018
$data2[18,3] = 1.028803158772578
$data2[18,4] = 1.035686457458766
$data2[18,5] = 1.011260454511378
$data2[19,0] = 1.027664993056495
$data2[19,1] = 1.021018435687277
$data2[19,2] = 1.026821841561943
$data2[19,3] = 1.027912569320712
$data2[19,4] = 1.034604285185606
$data2[19,5] = 1.010988859545932
$data2[20,0] = 1.027586509990163
$data2[20,1] = 1.020521858121801
$data2[20,2] = 1.026547147649343
$data2[20,3] = 1.027352892817891
$data2[20,4] = 1.033924544930363
$data2[20,5] = 1.010817902416904
$data2[21,0] = 1.027628606217047
$data2[21,1] = 1.020785101338793
$data2[21,2] = 1.026693504457306
$data2[21,3] = 1.027649562090435
$data2[21,4] = 1.034284825156833
$data2[21,5] = 1.010908548417755
$data2[22,0] = 1.027783436694908
$data2[22,1] = 1.021821508704644
$data2[22,2] = 1.027253345444725
$data2[22,3] = 1.028818094240572
$data2[22,4] = 1.035704611492333
$data2[22,5] = 1.011265004526705
$data2[23,0] = 1.027940447854658
$data2[23,1] = 1.023024291679433
$data2[23,2] = 1.027868987877461
$data2[23,3] = 1.030175296368112
$data2[23,4] = 1.037355091090547
$data2[23,5] = 1.011677799327959

$ws.Range("I2:N25").Value = $data2

Write-Output "Updated vm_pu values for 380 kV case"
